$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 28
$ws.Range("H28").Value = 630
$ws.Range("I28").Value = 278.6
$ws.Range("J28").Value = 1215.6666
$ws.Range("K28").Value = 278.6
$ws.Range("L28").Value = 1215.6666
$ws.Range("M28").Value = 206.4
$ws.Range("N28").Value = -2185.6666

# row 62
$ws.Range("H62").Value = 6833.0557
$ws.Range("I62").Value = 5773.1333
$ws.Range("J62").Value = 12132.667
$ws.Range("K62").Value = 5773.1333
$ws.Range("L62").Value = 12132.667
$ws.Range("M62").Value = -5149.1333
$ws.Range("N62").Value = -13380.667

# row 65
$ws.Range("H65").Value = 6833.0557
$ws.Range("I65").Value = 5773.1333
$ws.Range("J65").Value = 12132.667
$ws.Range("K65").Value = 28865.6665
$ws.Range("L65").Value = 60663.335
$ws.Range("M65").Value = -25745.6665
$ws.Range("N65").Value = -66903.33499999999

# row 129
$ws.Range("H129").Value = 2633660.5
$ws.Range("I129").Value = 33333868
$ws.Range("J129").Value = 2214.2856
$ws.Range("K129").Value = 100001604
$ws.Range("L129").Value = 6642.8568
$ws.Range("M129").Value = -99996604
$ws.Range("N129").Value = -16642.8568

# row 137
$ws.Range("H137").Value = 10000.333
$ws.Range("I137").Value = 10000.333
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 30000.999
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -27450.999

$ws = $wb.Worksheets.Item("ARM")
# row 97
$ws.Range("H97").Value = 366
$ws.Range("I97").Value = 282.5
$ws.Range("J97").Value = 700
$ws.Range("K97").Value = 282.5
$ws.Range("L97").Value = 700
$ws.Range("M97").Value = 213.5
$ws.Range("N97").Value = -1692

# row 110
$ws.Range("H110").Value = 546.5
$ws.Range("I110").Value = 546.5
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 546.5
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1498.5

# row 122
$ws.Range("H122").Value = 2204.1667
$ws.Range("I122").Value = 1404
$ws.Range("J122").Value = 3004.3333
$ws.Range("K122").Value = 4212
$ws.Range("L122").Value = 9012.999899999999
$ws.Range("M122").Value = -1762
$ws.Range("N122").Value = -13912.9999

# row 132
$ws.Range("H132").Value = 6396.143
$ws.Range("I132").Value = 3254.6
$ws.Range("J132").Value = 14250
$ws.Range("K132").Value = 9763.799999999999
$ws.Range("L132").Value = 42750
$ws.Range("M132").Value = -7233.799999999999
$ws.Range("N132").Value = -47810

# row 135
$ws.Range("H135").Value = 38999.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 38999.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 38999.5
$ws.Range("N135").Value = -49139.5

$ws = $wb.Worksheets.Item("BSM")
# row 92
$ws.Range("H92").Value = 19499.5
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 19499.5
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 19499.5
$ws.Range("N92").Value = -24491.5

# row 100
$ws.Range("H100").Value = 18663.334
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 18663.334
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 18663.334
$ws.Range("N100").Value = -20827.334

# row 134
$ws.Range("H134").Value = 4585.364
$ws.Range("I134").Value = 2711.125
$ws.Range("J134").Value = 9583.333000000001
$ws.Range("K134").Value = 8133.375
$ws.Range("L134").Value = 28749.999
$ws.Range("M134").Value = -5598.375
$ws.Range("N134").Value = -33819.999

$ws = $wb.Worksheets.Item("CRP")
# row 16
$ws.Range("H16").Value = 511
$ws.Range("I16").Value = 511
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 511
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -224
$ws.Range("N16").ClearContents()

# row 60
$ws.Range("H60").Value = 10000
$ws.Range("I60").Value = 10000
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 10000
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -9489

# row 105
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()

# row 113
$ws.Range("H113").Value = 511
$ws.Range("I113").Value = 511
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 511
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1659
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# row 40
$ws.Range("H40").Value = 30
$ws.Range("I40").Value = 30
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 120
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -51

# row 81
$ws.Range("H81").Value = 8000
$ws.Range("I81").Value = 8000
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 24000
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -22877

# row 84
$ws.Range("H84").Value = 8000
$ws.Range("I84").Value = 8000
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 72000
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -66384

# row 131
$ws.Range("H131").Value = 2749
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 2749
$ws.Range("K131").Value = 0
$ws.Range("L131").ClearContents()
$ws.Range("M131").Value = 8247
$ws.Range("N131").Value = -18327

# row 133
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = 0

# row 141
$ws.Range("H141").Value = 9950
$ws.Range("I141").Value = 9950
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 29850
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -24670

$ws = $wb.Worksheets.Item("GSM")
# row 101
$ws.Range("H101").Value = 1100597
$ws.Range("I101").Value = 1100597
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 1100597
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = -1097352
$ws.Range("N101").ClearContents()

# row 102
$ws.Range("H102").Value = 3376.1
$ws.Range("I102").Value = 2248.375
$ws.Range("J102").Value = 7887
$ws.Range("K102").Value = 2248.375
$ws.Range("L102").Value = 7887
$ws.Range("M102").Value = -626.375
$ws.Range("N102").Value = -11131

# row 132
$ws.Range("H132").Value = 6099.619
$ws.Range("I132").Value = 4505.8125
$ws.Range("J132").Value = 11199.8
$ws.Range("K132").Value = 13517.4375
$ws.Range("L132").Value = 33599.39999999999
$ws.Range("M132").Value = -10987.4375
$ws.Range("N132").Value = -38659.39999999999

# row 135
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("N135").Value = 0

$ws = $wb.Worksheets.Item("LTW")
# row 29
$ws.Range("H29").Value = 26066.666
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 26066.666
$ws.Range("K29").Value = 0
$ws.Range("L29").ClearContents()
$ws.Range("M29").Value = 26066.666
$ws.Range("N29").Value = -26656.666

# row 46
$ws.Range("H46").Value = 7700
$ws.Range("I46").Value = 7900
$ws.Range("J46").Value = 7500
$ws.Range("K46").Value = 7900
$ws.Range("L46").Value = 7500
$ws.Range("M46").Value = -7712
$ws.Range("N46").Value = -7876

$ws = $wb.Worksheets.Item("WVR")
# row 62
$ws.Range("H62").Value = 6833.3335
$ws.Range("I62").Value = 6500
$ws.Range("J62").Value = 7500
$ws.Range("K62").Value = 6500
$ws.Range("L62").Value = 7500
$ws.Range("M62").Value = -5876
$ws.Range("N62").Value = -8748

# row 65
$ws.Range("H65").Value = 6833.3335
$ws.Range("I65").Value = 6500
$ws.Range("J65").Value = 7500
$ws.Range("K65").Value = 32500
$ws.Range("L65").Value = 37500
$ws.Range("M65").Value = -29380
$ws.Range("N65").Value = -43740

# row 74
$ws.Range("H74").Value = 35000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 35000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 35000
$ws.Range("N74").Value = -36872

# row 77
$ws.Range("H77").Value = 35000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 35000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 105000
$ws.Range("N77").Value = -114360

# row 107
$ws.Range("H107").Value = 1457.75
$ws.Range("I107").Value = 1486.25
$ws.Range("J107").Value = 1400.75
$ws.Range("K107").Value = 4458.75
$ws.Range("L107").Value = 4202.25
$ws.Range("M107").Value = -2538.75
$ws.Range("N107").Value = -8042.25

# row 132
$ws.Range("H132").Value = 7965.3335
$ws.Range("I132").Value = 7965.3335
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 23896.0005
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -21366.0005
